# ------------------------------------------------------------------------
# Applies the commit's changes to the workbook:
#   1. Tidies up a now-superfluous explicit number-format override on
#      "Clients"!Q2:Q30 (the postal/box code column).
#   2. Inserts a blank spacer row into "Departments" (row 3), pushing the
#      existing "Home Nursing" row down to row 4.
#   3. Adds a new "Rebates-Purchases" worksheet (at the end of the tab
#      strip) containing a Rebates/Purchases report for five people across
#      twelve monthly periods.
#   4. Re-activates "Departments" so it remains the selected tab (adding
#      / populating the new sheet above leaves it activated instead).
# ------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Clients: drop the stray NumberFormat override on the code col ---
$clients = $wb.Worksheets.Item("Clients")
$clients.Range("Q2:Q30").Font.Color = 0

# --- 2. Departments: insert a blank row above the "Home Nursing" row ----
$dept = $wb.Worksheets.Item("Departments")
$dept.Rows("3:3").Insert()
$dept.Range("A3").Value = "      "

# --- 3. Add the Rebates-Purchases worksheet ------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Rebates-Purchases"

$ws.Columns.Item(1).ColumnWidth = 11.25

# Row 1: the 12 monthly period end-dates (23rd of each month, Jan-Dec 2023),
# each date spans two columns (one for Rebates, one for Purchases).
$periodDates = @(44949,44980,45008,45039,45069,45100,45130,45161,45192,45222,45253,45283)
$dateRow = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 12; $i++) {
    $dateRow[0, $i*2]     = $periodDates[$i]
    $dateRow[0, $i*2 + 1] = $periodDates[$i]
}
$ws.Range("C1:Z1").Value = $dateRow
$ws.Range("C1:Z1").NumberFormat = "d-mmm"

# Row 2: column headers.
$ws.Range("A2").Value = "Name*"
$ws.Range("B2").Value = "Group*"
$headerRow = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 12; $i++) {
    $headerRow[0, $i*2]     = "Rebates"
    $headerRow[0, $i*2 + 1] = "Purchases"
}
$ws.Range("C2:Z2").Value = $headerRow

# Rows 3-7: one person per row, Rebates/Purchases growing by 10 each period.
$people = @(
    @{ Name = "John Doe";       Group = "Group A" },
    @{ Name = "Jane Smith";     Group = "Group B" },
    @{ Name = "David Johnson";  Group = "Group C" },
    @{ Name = "Lisa Adams";     Group = "Group D" },
    @{ Name = "Mary Johnson";   Group = "Group E" }
)

for ($p = 0; $p -lt $people.Count; $p++) {
    $row = 3 + $p
    $ws.Range("A$row").Value = $people[$p].Name
    $ws.Range("B$row").Value = $people[$p].Group

    $base = ($p + 1) * 100
    $dataRow = New-Object 'object[,]' 1,24
    for ($i = 0; $i -lt 12; $i++) {
        $rebate = $base + $i * 10
        $dataRow[0, $i*2]     = $rebate
        $dataRow[0, $i*2 + 1] = $rebate * 10
    }
    $ws.Range("C$row`:Z$row").Value = $dataRow
}

$ws.Range("B2").Select()

# --- 4. Re-activate Departments so it's the sheet that's on top when ----
#        the file is saved (matching the saved tabSelected/selection).
$dept.Activate()
$dept.Range("D5").Select()
